$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.127006530761719
$ws.Range("B1").Value = 1.826100707054138
$ws.Range("D1").Value = 2.335582733154297
$ws.Range("E1").Value = 1.128461003303528
